# Automatische test-sync: 2025-08-05 17:03:50
# Appends a new log row (row 12) to the "Logs" sheet, mirroring the most
# recent "Kun jij dit even regelen?" entry with an updated timestamp, and
# bumps the matching "Planning / Afspraak" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 12
$logs.Cells.Item($newRow, 1).Value  = "Kun jij dit even regelen?"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #1: Kun jij dit even regelen?"
$logs.Cells.Item($newRow, 4).Value  = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-05 17:03:26"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Keep the dashboard tally for "Planning / Afspraak" (Dashboard!B2) in sync
# with the newly logged row.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = $dashboard.Range("B2").Value2 + 1

# Extend the conditional-formatting ranges so they keep covering the whole
# data area (previously rows 2-11, now rows 2-12) without disturbing the
# existing rules (same dxfId / operator / formula / priority per rule).
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$col`2:$col`11")
    $newRange = $logs.Range("$col`2:$col`12")
    $conditions = $oldRange.FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($newRange)
    }
}
